# Fig16.xlsx monthly snapshot roll-forward: January 2017 -> February 2017
# Commit message: "2017-02-13 snapshot - chunk 30"
#
# The workbook title / source caption strings move from "January 2017" to
# "February 2017", and the forecast-scenario data block (rows 85:110,
# columns B/C -- Distillate/Gasoline inventory levels for the new forecast
# path starting Nov-2016) is replaced with the refreshed STEO figures.
# Everything else (MIN/MAX helper columns E:H, the I:K range-width helper
# columns, the chart number caches, the external-link cache, and the
# calcChain) is formula-driven off these cells and recalculates
# automatically once Excel recalcs the workbook.
#
# The "Forecast" vertical marker line (scatter series anchored at
# A115:A116) also shifts one column to the right (x = 60 -> 61) since the
# forecast boundary itself moved forward a month.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fig16")

# --- Title (A2) and source caption (A111), both shared strings ---
$ws.Range("A2").Value = "Short-Term Energy Outlook, February 2017"
$ws.Range("A111").Value = "Source: Short-Term Energy Outlook, February 2017."

# --- Refreshed forecast data, rows 85-110 (Nov 2016 - Dec 2018) ---
$ws.Range("B85").Value = 160.173
$ws.Range("C85").Value = 233.416
$ws.Range("B86").Value = 164.07285714
$ws.Range("C86").Value = 236.88514285900001
$ws.Range("B87").Value = 170.17267548999999
$ws.Range("C87").Value = 258.57071535900002
$ws.Range("B88").Value = 162.34540000000001
$ws.Range("C88").Value = 250.27099999999999
$ws.Range("B89").Value = 157.9442
$ws.Range("C89").Value = 238.477
$ws.Range("B90").Value = 155.5179
$ws.Range("C90").Value = 232.00303
$ws.Range("B91").Value = 158.87520000000001
$ws.Range("C91").Value = 229.96593999999999
$ws.Range("B92").Value = 160.7234
$ws.Range("C92").Value = 230.46093999999999
$ws.Range("B93").Value = 166.15770000000001
$ws.Range("C93").Value = 230.27011999999999
$ws.Range("B94").Value = 169.48580000000001
$ws.Range("C94").Value = 225.99796000000001
$ws.Range("B95").Value = 167.3724
$ws.Range("C95").Value = 226.79655
$ws.Range("B96").Value = 160.4418
$ws.Range("C96").Value = 221.46883
$ws.Range("B97").Value = 162.33340000000001
$ws.Range("C97").Value = 229.90840000000003
$ws.Range("B98").Value = 167.21539999999999
$ws.Range("C98").Value = 241.36967999999999
$ws.Range("B99").Value = 163.8683
$ws.Range("C99").Value = 249.29425000000001
$ws.Range("B100").Value = 156.36500000000001
$ws.Range("C100").Value = 247.41799
$ws.Range("B101").Value = 152.19120000000001
$ws.Range("C101").Value = 239.28563
$ws.Range("B102").Value = 150.20830000000001
$ws.Range("C102").Value = 234.31788
$ws.Range("B103").Value = 153.9221
$ws.Range("C103").Value = 232.52828
$ws.Range("B104").Value = 155.93790000000001
$ws.Range("C104").Value = 232.90031000000002
$ws.Range("B105").Value = 161.5864
$ws.Range("C105").Value = 232.31128999999999
$ws.Range("B106").Value = 165.2927
$ws.Range("C106").Value = 228.36971
$ws.Range("B107").Value = 163.3938
$ws.Range("C107").Value = 228.60893000000002
$ws.Range("B108").Value = 156.57
$ws.Range("C108").Value = 223.45211
$ws.Range("B109").Value = 158.5341
$ws.Range("C109").Value = 231.86373999999998
$ws.Range("B110").Value = 163.608
$ws.Range("C110").Value = 243.87798000000001

# --- "Forecast" marker line x-position shifts one month forward ---
$ws.Range("A115").Value = 61
$ws.Range("A116").Value = 61

$excel.Calculate()
